$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "IsStartingFirstHalf" feature row is being removed/replaced: row 46 now documents
# IsStartingScndHalf (moved up from row 47) along with a note on why the
# "IsStartingFirstHalf" idea was dropped, highlighted with a red fill.
$ws.Range("A46").Value = "IsStartingScndHalf"
$ws.Range("B46").Value = "Removed could be unclear (turnovers on kick return/onside kicks)"
$ws.Range("A46:B46").Interior.Color = 192

# Row 47 used to hold "IsStartingScndHalf" - now blank (matches the blank A48 below it).
$ws.Range("A47").ClearContents()

# Remove the extra blank spacer row (old row 49), shifting everything below up by one.
$ws.Rows("49").Delete()

# Column B needs to be a bit wider to fit the new, longer note text.
$ws.Columns("B").ColumnWidth = 55.1667

# Reflect the scrolled viewport / active cell after the edits above.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 37
$ws.Range("A47").Select() | Out-Null
